# Apply the edits described by the diff.
$d = $word.ActiveDocument

# 1) "tretåig hackspett" -> "spillkråka och tretåig hackspett"
#    (occurs twice in the document; both occurrences get the same prefix)
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("tretåig hackspett", $true, $false, $false, $false, $false, `
               $true, 1, $false, "spillkråka och tretåig hackspett", 2, `
               $false, $false, $false, $false)

# 2) "Detta är en prioriterad art" -> "Dessa är prioriterade arter"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Detta är en prioriterad art", $true, $false, $false, $false, $false, `
               $true, 1, $false, "Dessa är prioriterade arter", 2, `
               $false, $false, $false, $false)

# 3) " denna art" -> " dessa arter"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(" denna art", $true, $false, $false, $false, $false, `
               $true, 1, $false, " dessa arter", 2, `
               $false, $false, $false, $false)

# 4) " arten" -> " arterna" (match whole word to avoid touching other words)
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result4 = $find.Execute("arten", $true, $true, $false, $false, $false, `
               $true, 1, $false, "arterna", 2, `
               $false, $false, $false, $false)
Write-Host "Step4 result: $result4"

# 5) Update the date "2026-02-10" -> "2026-02-11"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("2026-02-10", $true, $false, $false, $false, $false, `
               $true, 1, $false, "2026-02-11", 2, `
               $false, $false, $false, $false)
